# Update "Horarios Linea 141" workbook with the latest scrape results.
$wb = $excel.ActiveWorkbook

$newTime = "02:58:51"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 10"

$ws1.Range("A9").Value = $newTime
$ws1.Range("D9").Value = 0

$ws1.Range("A11").Value = $newTime
$ws1.Range("D11").Value = 50

$ws1.Range("A13").Value = $newTime
$ws1.Range("D13").Value = 63

$ws1.Range("A14").Value = $newTime
$ws1.Range("D14").Value = 107

$ws1.Range("A15").Value = $newTime
$ws1.Range("B15").Value = "04:53"
$ws1.Range("C15").Value = "11_ETCHEVERRY"
$ws1.Range("D15").Value = 115
$ws1.Range("E15").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"

$ws2.Range("A8").Value = $newTime
$ws2.Range("D8").Value = 0

$ws2.Range("A10").Value = $newTime
$ws2.Range("D10").Value = 107

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
